# Applies the Bahamut_Profits.xlsx profit-recalculation update (scheduled runner refresh).
# For each affected row, sets new H..N values; clears cells that the refresh removed (M82/M85);
# and adds cells that the refresh introduced (N94 on BSM).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 83336020
$ws.Range("H79").Value = 83336020
$ws.Range("H82").Value = 13000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 13000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 39000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -39812
$ws.Range("H85").Value = 13000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 13000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 39000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -41808
$ws.Range("H88").Value = 1325080.4
$ws.Range("I88").Value = 2735.7144
$ws.Range("J88").Value = 3176363
$ws.Range("K88").Value = 2735.7144
$ws.Range("L88").Value = 3176363
$ws.Range("M88").Value = -2329.7144
$ws.Range("N88").Value = -3177175
$ws.Range("H91").Value = 1325080.4
$ws.Range("I91").Value = 2735.7144
$ws.Range("J91").Value = 3176363
$ws.Range("K91").Value = 2735.7144
$ws.Range("L91").Value = 3176363
$ws.Range("M91").Value = -1331.7144
$ws.Range("N91").Value = -3179171
$ws.Range("H98").Value = 1575.3
$ws.Range("I98").Value = 1721.7646
$ws.Range("J98").Value = 745.3333
$ws.Range("K98").Value = 1721.7646
$ws.Range("L98").Value = 745.3333
$ws.Range("M98").Value = -223.7646
$ws.Range("N98").Value = -3741.3333
$ws.Range("H122").Value = 1575.3
$ws.Range("I122").Value = 1721.7646
$ws.Range("J122").Value = 745.3333
$ws.Range("K122").Value = 5165.293799999999
$ws.Range("L122").Value = 2235.9999
$ws.Range("M122").Value = -2715.293799999999
$ws.Range("N122").Value = -7135.9999
$ws.Range("H132").Value = 2262.7334
$ws.Range("I132").Value = 2377.9285
$ws.Range("K132").Value = 7133.7855
$ws.Range("M132").Value = -4603.7855
$ws.Range("H138").Value = 2827.19
$ws.Range("I138").Value = 943.119
$ws.Range("J138").Value = 4191.517
$ws.Range("K138").Value = 2829.357
$ws.Range("L138").Value = 12574.551
$ws.Range("M138").Value = 2310.643
$ws.Range("N138").Value = -22854.551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2037.25
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 2111.75
$ws.Range("K97").Value = 2000
$ws.Range("L97").Value = 2111.75
$ws.Range("M97").Value = -1504
$ws.Range("N97").Value = -3103.75
$ws.Range("H102").Value = 2634.4443
$ws.Range("I102").Value = 2686.1538
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2686.1538
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1064.1538
$ws.Range("N102").Value = -5744
$ws.Range("H124").Value = 24954.428
$ws.Range("J124").Value = 24954.428
$ws.Range("L124").Value = 24954.428
$ws.Range("N124").Value = -34774.428
$ws.Range("H125").Value = 20558.533
$ws.Range("J125").Value = 20558.533
$ws.Range("L125").Value = 20558.533
$ws.Range("N125").Value = -30398.533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 43480936
$ws.Range("I86").Value = 50002070
$ws.Range("J86").Value = 6700
$ws.Range("K86").Value = 50002070
$ws.Range("L86").Value = 6700
$ws.Range("M86").Value = -50000947
$ws.Range("N86").Value = -8946
$ws.Range("H89").Value = 43480936
$ws.Range("I89").Value = 50002070
$ws.Range("J89").Value = 6700
$ws.Range("K89").Value = 250010350
$ws.Range("L89").Value = 33500
$ws.Range("M89").Value = -250004734
$ws.Range("N89").Value = -44732
$ws.Range("H94").Value = 951.2727
$ws.Range("I94").Value = 975.1429000000001
$ws.Range("J94").Value = 450
$ws.Range("K94").Value = 975.1429000000001
$ws.Range("L94").Value = 450
$ws.Range("M94").Value = -524.1429000000001
$ws.Range("N94").Value = -1352
$ws.Range("H105").Value = 4787.875
$ws.Range("I105").Value = 4094
$ws.Range("J105").Value = 6473
$ws.Range("K105").Value = 4094
$ws.Range("L105").Value = 6473
$ws.Range("M105").Value = -2347
$ws.Range("N105").Value = -9967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 7885.75
$ws.Range("J94").Value = 8361.700000000001
$ws.Range("L94").Value = 8361.700000000001
$ws.Range("N94").Value = -9263.700000000001
$ws.Range("H112").Value = 21943.334
$ws.Range("J112").Value = 21943.334
$ws.Range("L112").Value = 21943.334
$ws.Range("N112").Value = -24897.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1283.7858
$ws.Range("I122").Value = 431.33334
$ws.Range("J122").Value = 1516.2727
$ws.Range("K122").Value = 3882.00006
$ws.Range("L122").Value = 13646.4543
$ws.Range("M122").Value = -1432.00006
$ws.Range("N122").Value = -18546.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1112.3529
$ws.Range("I97").Value = 1283.3334
$ws.Range("J97").Value = 702
$ws.Range("K97").Value = 1283.3334
$ws.Range("L97").Value = 702
$ws.Range("M97").Value = -787.3334
$ws.Range("N97").Value = -1694
$ws.Range("H102").Value = 2146.5715
$ws.Range("I102").Value = 2070.6667
$ws.Range("J102").Value = 2203.5
$ws.Range("K102").Value = 2070.6667
$ws.Range("L102").Value = 2203.5
$ws.Range("M102").Value = -448.6667000000002
$ws.Range("N102").Value = -5447.5
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
$ws.Range("H113").Value = 3902.4
$ws.Range("I113").Value = 4566.5
$ws.Range("J113").Value = 2906.25
$ws.Range("K113").Value = 4566.5
$ws.Range("L113").Value = 2906.25
$ws.Range("M113").Value = -2396.5
$ws.Range("N113").Value = -7246.25
$ws.Range("H126").Value = 1934.68
$ws.Range("I126").Value = 1946.3914
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 5839.174199999999
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -3369.174199999999
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 12347857
$ws.Range("I100").Value = 18520618
$ws.Range("J100").Value = 2334.6667
$ws.Range("K100").Value = 18520618
$ws.Range("L100").Value = 2334.6667
$ws.Range("M100").Value = -18520077
$ws.Range("N100").Value = -3416.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 39000
$ws.Range("J16").Value = 39000
$ws.Range("L16").Value = 39000
$ws.Range("N16").Value = -39584

Write-Host "Bahamut_Profits.xlsx scheduled-runner update applied."
